$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; existing rows 10-20 shift down to 11-21.
$ws.Rows("10").Insert()

# Populate the newly inserted row 10 with the new weekly price record.
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = "Vega Monumental Concepción"
$ws.Range("C10").Value = "Bíobío"
$ws.Range("D10").Value = 45014
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100101
$ws.Range("H10").Value = "Berries"
$ws.Range("I10").Value = 100101004
$ws.Range("J10").Value = "Frambuesa"
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 7000
$ws.Range("O10").Value = 7500
$ws.Range("P10").Value = 7250
$ws.Range("Q10").Value = "`$/bandeja 2 kilos"
$ws.Range("R10").Value = "Región de Ñuble"
$ws.Range("S10").Value = 3625
$ws.Range("T10").Value = 2
